$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Dolichotis / erythyro.2 component) - fix "erythryo" -> "erythyro" typo
# and correct values so the locus columns hold the proper labels.
$ws.Range("A2").Value = "erythyro.2-Dolichotis_patagonum"
$ws.Range("B2").Value = "Dolichotis_patagonum"
$ws.Range("C2").Value = "erythyro.2-dolichotis"

# Row 3 (Indri / erythyro.1 component) - same typo fix.
$ws.Range("A3").Value = "erythyro.1-Indri_indri"
$ws.Range("B3").Value = "Indri_indri"
$ws.Range("C3").Value = "erythyro.1-indri"

$ws.Range("E3").Value = "erythyro.1"
$ws.Range("E2").Value = "erythyro.2"

# Update the selected cell in the sheet view.
$ws.Range("B10").Select()
